$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1) - copy formatting from the adjacent "sum"
# header (G1) so the bold/centered/bordered header style carries over,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data cell (H2) with numeric value 0.
$ws.Range("H2").Value = 0
